$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "63.838.13"
$ws.Cells.Item(2, 5).Value = "  -1.63%  "
$ws.Cells.Item(3, 4).Value = "3.053.22"
$ws.Cells.Item(3, 5).Value = "  -1.94%  "
$ws.Cells.Item(4, 5).Value = "  -0.05%  "
$ws.Cells.Item(5, 4).Value = "'559.09"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -0.10%  "
$ws.Cells.Item(6, 4).Value = "'142.52"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -1.95%  "
$ws.Cells.Item(7, 5).Value = "  +0.06%  "
$ws.Cells.Item(8, 4).Value = "3.052.38"
$ws.Cells.Item(8, 5).Value = "  -1.91%  "
$ws.Cells.Item(9, 4).Value = "'0.514"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +2.57%  "
$ws.Cells.Item(10, 4).Value = "'0.154"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +0.41%  "
$ws.Cells.Item(11, 4).Value = "'6.24"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -11.81%  "
$ws.Cells.Item(12, 4).Value = "'0.494"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = "  +6.84%  "
$ws.Cells.Item(13, 4).Value = "'0.0000230"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = "  +1.06%  "
$ws.Cells.Item(14, 4).Value = "'35.69"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = "  +0.42%  "
$ws.Cells.Item(15, 4).Value = "3.552.03"
$ws.Cells.Item(15, 5).Value = "  -1.33%  "
$ws.Cells.Item(16, 4).Value = "63.885.05"
$ws.Cells.Item(16, 5).Value = "  -1.61%  "
$ws.Cells.Item(17, 4).Value = "3.050.24"
$ws.Cells.Item(17, 5).Value = "  -2.01%  "
$ws.Cells.Item(18, 5).Value = "  +0.47%  "
$ws.Cells.Item(19, 4).Value = "'6.79"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +0.39%  "
$ws.Cells.Item(20, 4).Value = "'475.57"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -2.68%  "
$ws.Cells.Item(21, 4).Value = "'14.04"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +1.60%  "
$ws.Cells.Item(22, 4).Value = "'0.684"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +1.22%  "
$ws.Cells.Item(23, 4).Value = "'14.59"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +9.72%  "
$ws.Cells.Item(24, 4).Value = "'7.55"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -0.16%  "
$ws.Cells.Item(25, 4).Value = "'82.57"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +1.73%  "
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -0.60%  "
$ws.Cells.Item(27, 5).Value = "  -0.91%  "
$ws.Cells.Item(28, 4).Value = "'8.14"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  +1.58%  "
$ws.Cells.Item(29, 4).Value = "'2.04"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -1.25%  "
$ws.Cells.Item(30, 4).Value = "'0.999"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -0.09%  "
$ws.Cells.Item(31, 4).Value = "'26.26"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +0.03%  "
$ws.Cells.Item(32, 5).Value = "  -1.56%  "
$ws.Cells.Item(33, 5).Value = "  -0.52%  "
$ws.Cells.Item(34, 4).Value = "'5.77"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -0.21%  "
$ws.Cells.Item(35, 4).Value = "'6.22"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = "  +1.12%  "
$ws.Cells.Item(36, 4).Value = "'54.51"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -1.97%  "
$ws.Cells.Item(37, 4).Value = "'0.0410"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +0.35%  "
$ws.Cells.Item(38, 4).Value = "'447.08"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(39, 4).Value = "'0.0814"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -1.78%  "
$ws.Cells.Item(40, 4).Value = "'2.83"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +4.25%  "
$ws.Cells.Item(41, 4).Value = "3.010.98"
$ws.Cells.Item(41, 5).Value = "  -0.99%  "
$ws.Cells.Item(42, 5).Value = "  -0.77%  "
$ws.Cells.Item(43, 4).Value = "'8.28"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -0.61%  "
$ws.Cells.Item(44, 4).Value = "'0.269"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +2.80%  "
$ws.Cells.Item(45, 4).Value = "'28.15"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -0.73%  "
$ws.Cells.Item(46, 4).Value = "'2.26"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +7.96%  "
$ws.Cells.Item(48, 5).Value = "  +0.85%  "
$ws.Cells.Item(49, 4).Value = "'117.75"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -0.42%  "
$ws.Cells.Item(50, 4).Value = "0.0₃0514"
$ws.Cells.Item(50, 5).Value = "  -1.37%  "
$ws.Cells.Item(51, 5).Value = "  +0.62%  "
